$wb = $excel.ActiveWorkbook

# 1. Target sheet: assembly value "GRCh38.p7" -> "GRCh38" (D2)
$target = $wb.Worksheets.Item("Target")
$target.Cells.Item(2, 4).Value = "GRCh38"

# 2. ExperimentLayout sheet: cell_pool column (E) values change from
#    "<prefix>.<n>" (e.g. 82.3) to just "<n>" (e.g. 3) - the leading
#    plate-prefix digits are dropped, keeping only the trailing digit.
$layout = $wb.Worksheets.Item("ExperimentLayout")

# Each entry is encoded as "row:newValue"
$poolUpdates = @(
    "3:3", "4:3", "5:3", "6:3", "7:1", "8:1", "9:1", "10:1", "11:2", "13:2", "14:1", "15:3", "16:3",
    "17:3", "18:3", "19:1", "20:1", "21:1", "22:1", "23:2", "25:2", "26:1", "27:3", "28:3", "29:3",
    "30:3", "31:1", "32:1", "33:1", "34:1", "35:2", "37:2", "38:1", "39:3", "40:3", "41:3", "42:3",
    "43:1", "44:1", "45:1", "46:1", "47:2", "49:2", "50:1", "51:3", "52:3", "53:3", "54:3", "55:1",
    "56:1", "57:1", "58:1", "59:2", "61:2", "62:1", "63:3", "64:3", "65:3", "66:3", "67:1", "68:1",
    "69:1", "70:1", "71:2", "73:2", "74:1", "75:3", "76:3", "77:3", "78:3", "79:1", "80:1", "81:1",
    "82:2", "83:2", "85:2", "86:3", "87:3", "88:3", "89:3", "90:1", "91:1", "92:1", "93:1", "94:2",
    "95:2", "99:2", "100:2", "102:2", "103:2", "104:2", "106:2", "107:2", "108:2", "110:2", "111:2",
    "112:2", "113:2", "114:2", "115:2", "116:2", "117:2", "119:2", "120:2", "122:2", "123:2"
)

foreach ($entry in $poolUpdates) {
    $parts = $entry.Split(":")
    $row = [int]$parts[0]
    $newValue = [int]$parts[1]
    $layout.Cells.Item($row, 5).Value = $newValue
}

Write-Output "Done: updated assembly cell and $($poolUpdates.Count) cell_pool values."